$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4281
$ws.Range("I40").Value = 3888.3333
$ws.Range("K40").Value = 3888.3333
$ws.Range("M40").Value = -3713.3333

$ws.Range("H62").Value = 3141.125
$ws.Range("I62").Value = 2693.75
$ws.Range("J62").Value = 3588.5
$ws.Range("K62").Value = 2693.75
$ws.Range("L62").Value = 3588.5
$ws.Range("M62").Value = -2069.75
$ws.Range("N62").Value = -4836.5

$ws.Range("H65").Value = 3141.125
$ws.Range("I65").Value = 2693.75
$ws.Range("J65").Value = 3588.5
$ws.Range("K65").Value = 13468.75
$ws.Range("L65").Value = 17942.5
$ws.Range("M65").Value = -10348.75
$ws.Range("N65").Value = -24182.5

$ws.Range("H74").Value = 1126959.8
$ws.Range("I74").Value = 1126959.8
$ws.Range("K74").Value = 1126959.8
$ws.Range("M74").Value = -1126023.8

$ws.Range("H77").Value = 1126959.8
$ws.Range("I77").Value = 1126959.8
$ws.Range("K77").Value = 5634799
$ws.Range("M77").Value = -5630119

$ws.Range("H111").Value = 3548.0833
$ws.Range("I111").Value = 3357.6667
$ws.Range("K111").Value = 10073.0001
$ws.Range("M111").Value = -7006.000100000001

$ws.Range("H132").Value = 3471.0667
$ws.Range("I132").Value = 3361.8572
$ws.Range("K132").Value = 10085.5716
$ws.Range("M132").Value = -7555.571599999999

$ws.Range("H137").Value = 3797.8667
$ws.Range("I137").Value = 3387
$ws.Range("K137").Value = 10161
$ws.Range("M137").Value = -7611

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2875.1353
$ws.Range("I32").Value = 2864.2058
$ws.Range("K32").Value = 2864.2058
$ws.Range("M32").Value = -2577.2058

$ws.Range("H102").Value = 14288101
$ws.Range("I102").Value = 20002142
$ws.Range("K102").Value = 20002142
$ws.Range("M102").Value = -20000520

$ws.Range("H122").Value = 3893
$ws.Range("I122").Value = 3081.1924
$ws.Range("K122").Value = 9243.5772
$ws.Range("M122").Value = -6793.5772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H38").Value = 3336
$ws.Range("J38").Value = 3336
$ws.Range("L38").Value = 3336
$ws.Range("N38").Value = -4168

$ws.Range("H86").Value = 4190.222
$ws.Range("J86").Value = 4184.5
$ws.Range("L86").Value = 4184.5
$ws.Range("N86").Value = -6430.5

$ws.Range("H89").Value = 4190.222
$ws.Range("J89").Value = 4184.5
$ws.Range("L89").Value = 20922.5
$ws.Range("N89").Value = -32154.5

$ws.Range("H103").Value = 23663.334
$ws.Range("J103").Value = 23663.334
$ws.Range("L103").Value = 23663.334
$ws.Range("N103").Value = -26007.334

$ws.Range("H105").Value = 2108.625
$ws.Range("I105").Value = 2081.2856
$ws.Range("K105").Value = 2081.2856
$ws.Range("M105").Value = -334.2856000000002

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 12616.667
$ws.Range("I31").Value = 8675
$ws.Range("J31").Value = 20500
$ws.Range("K31").Value = 8675
$ws.Range("L31").Value = 20500
$ws.Range("M31").Value = -8380
$ws.Range("N31").Value = -21090

$ws.Range("H34").Value = 12616.667
$ws.Range("I34").Value = 8675
$ws.Range("J34").Value = 20500
$ws.Range("K34").Value = 8675
$ws.Range("L34").Value = 20500
$ws.Range("M34").Value = -8473
$ws.Range("N34").Value = -20904

$ws.Range("H95").Value = 14262
$ws.Range("J95").Value = 14262
$ws.Range("L95").Value = 14262
$ws.Range("N95").Value = -19754

$ws.Range("H96").Value = 17199.285
$ws.Range("J96").Value = 17199.285
$ws.Range("L96").Value = 17199.285
$ws.Range("N96").Value = -22691.285

$ws.Range("H105").Value = 5630962.5
$ws.Range("I105").Value = 5630962.5
$ws.Range("K105").Value = 5630962.5
$ws.Range("M105").Value = -5629215.5

$ws.Range("H107").Value = 112470
$ws.Range("I107").Value = 763.6
$ws.Range("J107").Value = 252103
$ws.Range("K107").Value = 763.6
$ws.Range("L107").Value = 252103
$ws.Range("M107").Value = 1156.4
$ws.Range("N107").Value = -255943

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H7").Value = 5000247
$ws.Range("J7").Value = 495
$ws.Range("L7").Value = 1485
$ws.Range("N7").Value = -1709

$ws.Range("H18").Value = 5000
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("M18").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2942.7368
$ws.Range("I102").Value = 3053.647
$ws.Range("K102").Value = 3053.647
$ws.Range("M102").Value = -1431.647

$ws.Range("H126").Value = 3000
$ws.Range("I126").Value = 3000
$ws.Range("J126").Value = 3000
$ws.Range("K126").Value = 9000
$ws.Range("L126").Value = 9000
$ws.Range("M126").Value = -6530
$ws.Range("N126").Value = -13940

$ws.Range("H132").Value = 3293574.5
$ws.Range("I132").Value = 3680620
$ws.Range("K132").Value = 11041860
$ws.Range("M132").Value = -11039330

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2436.2104
$ws.Range("I46").Value = 1772.7333
$ws.Range("K46").Value = 1772.7333
$ws.Range("M46").Value = -1584.7333

$ws.Range("H122").Value = 3381.75
$ws.Range("I122").Value = 3314.9614
$ws.Range("K122").Value = 9944.8842
$ws.Range("M122").Value = -7494.8842

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 8998.4
$ws.Range("I62").Value = 12163
$ws.Range("J62").Value = 7642.143
$ws.Range("K62").Value = 12163
$ws.Range("L62").Value = 7642.143
$ws.Range("M62").Value = -11539
$ws.Range("N62").Value = -8890.143

$ws.Range("H65").Value = 8998.4
$ws.Range("I65").Value = 12163
$ws.Range("J65").Value = 7642.143
$ws.Range("K65").Value = 60815
$ws.Range("L65").Value = 38210.715
$ws.Range("M65").Value = -57695
$ws.Range("N65").Value = -44450.715

$ws.Range("H68").Value = 15000
$ws.Range("J68").Value = 15000
$ws.Range("L68").Value = 15000
$ws.Range("N68").Value = -16622

$ws.Range("H71").Value = 15000
$ws.Range("J71").Value = 15000
$ws.Range("L71").Value = 45000
$ws.Range("N71").Value = -53112

$ws.Range("H81").Value = 3342.8333
$ws.Range("I81").Value = 3342.8333
$ws.Range("J81").Value = 0
$ws.Range("K81").Value = 6685.6666
$ws.Range("L81").Value = 0
$ws.Range("M81").Value = -5624.6666
$ws.Range("N81").Value = $null

$ws.Range("H84").Value = 3342.8333
$ws.Range("I84").Value = 3342.8333
$ws.Range("J84").Value = 0
$ws.Range("K84").Value = 33428.333
$ws.Range("L84").Value = 0
$ws.Range("M84").Value = -28124.333
$ws.Range("N84").Value = $null

$ws.Range("H136").Value = 45457656
$ws.Range("I136").Value = 45457656
$ws.Range("K136").Value = 136372968
$ws.Range("M136").Value = -136370418

Write-Output "Spriggan_Profits update applied"
